$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45178 = 2023-09-09) for
# every data row (rows 2-90). The update bumps that date forward by one day
# (45179 = 2023-09-10) for all of those rows.
$ws.Range("C2:C90").Value = 45179
